$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParaIndex($substr) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

function Replace-Paragraphs($firstSubstr, $lastSubstr, $xml) {
    $iFirst = Find-ParaIndex($firstSubstr)
    $iLast = Find-ParaIndex($lastSubstr)
    $start = $d.Paragraphs($iFirst).Range.Start
    $end = $d.Paragraphs($iLast).Range.End
    $rng = $d.Range($start, $end)
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "propic" default path: split "C:\Users\Velios\Desktop\Uni\cd" so that
#    "Velios" is wrapped in spellcheck proofErr markers.
# ---------------------------------------------------------------------------
$xml1 = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:tab/></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>propic</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>VARCHAR(</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>1024)</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> DEFAULT &#8216;</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>C:\Users\</w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Velios</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>\Desktop\Uni\cd</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>\</w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>defaultpic</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>&#8217;</w:t></w:r>" +
  "</w:p>"
Replace-Paragraphs "Velios" "Velios" $xml1

# ---------------------------------------------------------------------------
# 2) "fav BOOLEAN NOT NULL" -> "favorite BOOLEAN NOT NULL DEFAULT 'false'"
#    (DEFAULT 'false' is colored red)
# ---------------------------------------------------------------------------
$xml2 = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:tab/><w:t>fav</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>orite</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> BOOLEAN</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> NOT NULL</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:r><w:rPr><w:color w:val=`"FF0000`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>DEFAULT &#8216;false&#8217;</w:t></w:r>" +
  "</w:p>"
Replace-Paragraphs "fav BOOLEAN" "fav BOOLEAN" $xml2

# ---------------------------------------------------------------------------
# 3) "welcome " -> "status " (new column name)
# ---------------------------------------------------------------------------
$xml3 = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>status</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>VARCHAR(</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">1024) </w:t></w:r>" +
  "</w:p>"
Replace-Paragraphs "welcome" "welcome" $xml3

# ---------------------------------------------------------------------------
# 4) email CHECK constraint: ' LIKE '_%_%.__%')' -> ' LIKE '_%@_%.__%')'
# ---------------------------------------------------------------------------
$xml4 = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">ADD CONSTRAINT </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>controllo_email</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> CHECK </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>( email</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> LIKE &#8216;_%</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>@_%</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>.__%&#8217;)</w:t></w:r>" +
  "</w:p>"
Replace-Paragraphs "controllo_email" "controllo_email" $xml4

# ---------------------------------------------------------------------------
# 5) & 6) password / nickname CHECK constraints: rewrite using LENGTH(), and
#    drop the stray blank paragraph + duplicate "ALTER TABLE User" paragraph
#    that used to sit between the two constraints.
# ---------------------------------------------------------------------------
$xmlPw = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">ADD CONSTRAINT </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>controllo_pw</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> CHECK (</w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>LENGTH(</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>password</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>)&gt;7</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>)</w:t></w:r>" +
  "</w:p>"
$xmlNick = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">ADD CONSTRAINT </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>controllo_nick</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> CHECK (</w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>LENGTH(</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>nickname)&gt;2</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>)</w:t></w:r>" +
  "</w:p>"
Replace-Paragraphs "controllo_pw" "controllo_nick" ($xmlPw + $xmlNick)
